$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data uses "Vijay Shankar" followed by a NON-BREAKING SPACE
# (U+00A0), matching the existing rows 2-6 in column F.
$batsman = "Vijay Shankar" + [char]0x00A0

# New rows (7-11) to append after the existing data (rows 1-6).
# Values mirror earlier rows but re-ordered/duplicated as per source diff.
$rows = @(
    @(" Dubai (DSC)", " October 13 2020", "Super Kings won by 20 runs", "Sunrisers Hyderabad", "Chennai Super Kings", $batsman, "12", "7", "0", "1", "171.42"),
    @(" Dubai (DSC)", " October 22 2020", "Sunrisers won by 8 wickets (with 11 balls remaining)", "Sunrisers Hyderabad", "Rajasthan Royals", $batsman, "52", "51", "6", "0", "101.96"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Sunrisers Hyderabad", "Kolkata Knight Riders", $batsman, "7", "10", "0", "0", "70.00"),
    @(" Dubai (DSC)", " September 21 2020", "RCB won by 10 runs", "Sunrisers Hyderabad", "Royal Challengers Bangalore", $batsman, "0", "1", "0", "0", "0.00"),
    @(" Dubai (DSC)", " October 24 2020", "Kings XI won by 12 runs", "Sunrisers Hyderabad", "Kings XI Punjab", $batsman, "26", "27", "4", "0", "96.29")
)

$startRow = 7

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowNum = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($rowNum, $c + 1)
        # Prefix with an apostrophe so numeric-looking strings (e.g. "12",
        # "171.42") stay text, matching the existing rows (t="str"),
        # without touching the cell's number format/style.
        $cell.Value = "'" + $rowData[$c]
        $cell.Style = "Normal"
    }
}

$endRow = $startRow + $rows.Count - 1

# The whole populated range (A1:K11) holds numeric-looking values stored
# as text -- mark the "number stored as text" warning ignored over the
# full extent, matching how Excel widens the ignoredErrors sqref as the
# table grows.
$full = $ws.Range("A1:K$endRow")
try {
    $full.Errors.Item(3).Ignore = $true
} catch {
}
